# The document contains six "<id>p088r_aN</id>" markers (N = 1..6), each
# split across three runs: a Courier-New-styled "<id>" run, a plain run
# holding "p088r_aN", and a Courier-New-styled "</id>" run. The edit
# collapses each trio into a single run reading "<id>p088r_N</id>" (the
# "a" is dropped), keeping the Courier-New / color 7f6000 / size 18
# formatting of the surrounding "<id>"/"</id>" runs.
#
# A plain Find & Replace across the whole "<id>p088r_aN</id>" span makes
# Word merge the three runs into one run that carries the formatting of
# the first run in the span (the "<id>" run) - exactly the target state.

$d = $word.ActiveDocument

for ($i = 1; $i -le 6; $i++) {
    $old = "<id>p088r_a$i</id>"
    $new = "<id>p088r_$i</id>"
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "WARNING: pattern not found for index $i : $old"
    }
}
